# Daily Status Tracker - latest status update
# (1) Overview: a couple of the rolled-up counters moved, and the "active"
#     tab/selection moved from Overview to Details.
# (2) Details: several rows' Testing/Development status + ETA dates were
#     refreshed to reflect progress made since the last update.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsDetails  = $wb.Worksheets.Item("Details")

# ---------------------------------------------------------------------
# Overview sheet: bump the rolled-up Development/Testing counters
# ---------------------------------------------------------------------
$wsOverview.Range("D4").Value = 27
$wsOverview.Range("E4").Value = 17

# ---------------------------------------------------------------------
# Details sheet: row 22 - Testing finished ("In Progress" -> "Done")
# ---------------------------------------------------------------------
$wsDetails.Range("H22").Value = "Done"

# Rows 27-30: ETA slipped a few days (2015-12-11 -> 2015-12-15)
$wsDetails.Range("J27").Value = 42353
$wsDetails.Range("J28").Value = 42353
$wsDetails.Range("J29").Value = 42353
$wsDetails.Range("J30").Value = 42353

# Row 31: Testing cell now center-aligned (matches the rest of the column)
$wsDetails.Range("H31").HorizontalAlignment = -4108   # xlCenter
$wsDetails.Range("J31").Value = 42353

# Row 32: ETA update only
$wsDetails.Range("J32").Value = 42353

# Row 33: Testing cell now center-aligned (matches the rest of the column)
$wsDetails.Range("H33").HorizontalAlignment = -4108   # xlCenter
$wsDetails.Range("J33").Value = 42353

# Rows 34-36: Development finished ("In Progress" -> "Done"), Testing now
# "In Progress", plus the ETA update
$wsDetails.Range("G34").Value = "Done"
$wsDetails.Range("H34").Value = "In Progress"
$wsDetails.Range("H34").HorizontalAlignment = -4108   # xlCenter
$wsDetails.Range("J34").Value = 42353

$wsDetails.Range("G35").Value = "Done"
$wsDetails.Range("H35").Value = "In Progress"
$wsDetails.Range("H35").HorizontalAlignment = -4108   # xlCenter
$wsDetails.Range("J35").Value = 42353

$wsDetails.Range("G36").Value = "Done"
$wsDetails.Range("H36").Value = "In Progress"
$wsDetails.Range("H36").HorizontalAlignment = -4108   # xlCenter
$wsDetails.Range("J36").Value = 42353

# ---------------------------------------------------------------------
# View state: Details is now the active tab, scrolled down a bit with a
# new selection; Overview keeps a plain (non-active) selection.
# ---------------------------------------------------------------------
$wsOverview.Activate()
$wsOverview.Range("F4").Select()

$wsDetails.Activate()
$winDetails = $excel.ActiveWindow
$winDetails.ScrollRow = 8
$winDetails.ScrollColumn = 2
$wsDetails.Range("G37").Select()
